$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row.
# D-column cells are forced to Text format before assignment so that
# values such as "304.91" or "39.501.09" are preserved exactly as
# text (matching the original inline-string / text cell content)
# instead of being auto-converted into numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.501.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.296.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.70%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.91"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "81.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.77%  "

$ws.Range("E7").Value = "  -4.12%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.468"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0783"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "28.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.46%  "

$ws.Range("E12").Value = "  -0.63%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.649.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.304.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.730"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.71%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "39.423.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0876"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "66.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.14%  "

$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.06%  "

$ws.Range("E26").Value = "  -4.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.58%  "

$ws.Range("E28").Value = "  -0.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.51%  "

$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0698"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.11%  "

$ws.Range("E36").Value = "  -1.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0960"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.76%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.949.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.69%  "

$ws.Range("E44").Value = "  -6.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.93%  "

$ws.Range("E47").Value = "  -9.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.515.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "89.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "48.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.85%  "

# Row 30/31 swap: InjectiveProtocol moves above Monero, with updated values
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.65%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "149.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.82%  "

